$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) of the last existing data row (51) down to the new rows (52-60)
$ws.Range("A51:V51").Copy()
$ws.Range("A52:V60").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 52
$ws.Range("B52").Value = "spain"
$ws.Range("C52").Value = "laliga"
$ws.Range("D52").Value = "2023-2024"
$ws.Range("A52").Value = 51
$ws.Range("E52").Value = 45192.58333333334
$ws.Range("F52").Value = "Girona"
$ws.Range("G52").Value = 5
$ws.Range("H52").Value = "Mallorca"
$ws.Range("I52").Value = 3
$ws.Range("J52").Value = 1.89
$ws.Range("K52").Value = "05/09/2023 12:02"
$ws.Range("L52").Value = 1.91
$ws.Range("M52").Value = "23/09/2023 13:55"
$ws.Range("N52").Value = 3.34
$ws.Range("O52").Value = "05/09/2023 12:02"
$ws.Range("P52").Value = 3.4
$ws.Range("Q52").Value = "23/09/2023 13:58"
$ws.Range("R52").Value = 4.33
$ws.Range("S52").Value = "05/09/2023 12:02"
$ws.Range("T52").Value = 4.88
$ws.Range("U52").Value = "23/09/2023 13:55"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/spain/laliga/girona-mallorca/0GnjIxkh/"

# Row 53
$ws.Range("B53").Value = "spain"
$ws.Range("C53").Value = "laliga"
$ws.Range("D53").Value = "2023-2024"
$ws.Range("A53").Value = 52
$ws.Range("E53").Value = 45192.67708333334
$ws.Range("F53").Value = "Osasuna"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = "Sevilla"
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 2.41
$ws.Range("K53").Value = "05/09/2023 12:02"
$ws.Range("L53").Value = 2.45
$ws.Range("M53").Value = "23/09/2023 16:14"
$ws.Range("N53").Value = 3.16
$ws.Range("O53").Value = "05/09/2023 12:02"
$ws.Range("P53").Value = 3.24
$ws.Range("Q53").Value = "23/09/2023 16:13"
$ws.Range("R53").Value = 3.09
$ws.Range("S53").Value = "05/09/2023 12:02"
$ws.Range("T53").Value = 3.28
$ws.Range("U53").Value = "23/09/2023 16:14"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/spain/laliga/osasuna-sevilla/CjjfHd4b/"

# Row 54
$ws.Range("B54").Value = "spain"
$ws.Range("C54").Value = "laliga"
$ws.Range("D54").Value = "2023-2024"
$ws.Range("A54").Value = 53
$ws.Range("E54").Value = 45192.77083333334
$ws.Range("F54").Value = "Barcelona"
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = "Celta Vigo"
$ws.Range("I54").Value = 2
$ws.Range("J54").Value = 1.42
$ws.Range("K54").Value = "05/09/2023 12:02"
$ws.Range("L54").Value = 1.26
$ws.Range("M54").Value = "23/09/2023 18:08"
$ws.Range("N54").Value = 4.91
$ws.Range("O54").Value = "05/09/2023 12:02"
$ws.Range("P54").Value = 6.77
$ws.Range("Q54").Value = "23/09/2023 18:29"
$ws.Range("R54").Value = 7.7
$ws.Range("S54").Value = "05/09/2023 12:02"
$ws.Range("T54").Value = 12.3
$ws.Range("U54").Value = "23/09/2023 18:29"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/spain/laliga/barcelona-celta-vigo/v1bsKbKu/"

# Row 55
$ws.Range("B55").Value = "spain"
$ws.Range("C55").Value = "laliga"
$ws.Range("D55").Value = "2023-2024"
$ws.Range("A55").Value = 54
$ws.Range("E55").Value = 45192.875
$ws.Range("F55").Value = "Almeria"
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = "Valencia"
$ws.Range("I55").Value = 2
$ws.Range("J55").Value = 2.82
$ws.Range("K55").Value = "05/09/2023 12:02"
$ws.Range("L55").Value = 2.72
$ws.Range("M55").Value = "23/09/2023 20:59"
$ws.Range("N55").Value = 3.34
$ws.Range("O55").Value = "05/09/2023 12:02"
$ws.Range("P55").Value = 3.43
$ws.Range("Q55").Value = "23/09/2023 20:59"
$ws.Range("R55").Value = 2.63
$ws.Range("S55").Value = "05/09/2023 12:02"
$ws.Range("T55").Value = 2.76
$ws.Range("U55").Value = "23/09/2023 20:59"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/spain/laliga/almeria-valencia/ELsWCsDC/"

# Row 56
$ws.Range("B56").Value = "spain"
$ws.Range("C56").Value = "laliga"
$ws.Range("D56").Value = "2023-2024"
$ws.Range("A56").Value = 55
$ws.Range("E56").Value = 45193.58333333334
$ws.Range("F56").Value = "Real Sociedad"
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = "Getafe"
$ws.Range("I56").Value = 3
$ws.Range("J56").Value = 1.56
$ws.Range("K56").Value = "05/09/2023 12:02"
$ws.Range("L56").Value = 1.71
$ws.Range("M56").Value = "24/09/2023 13:40"
$ws.Range("N56").Value = 3.63
$ws.Range("O56").Value = "05/09/2023 12:02"
$ws.Range("P56").Value = 3.43
$ws.Range("Q56").Value = "24/09/2023 13:40"
$ws.Range("R56").Value = 6.83
$ws.Range("S56").Value = "05/09/2023 12:02"
$ws.Range("T56").Value = 6.71
$ws.Range("U56").Value = "24/09/2023 13:40"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/spain/laliga/real-sociedad-getafe/2uq7EflH/"

# Row 57
$ws.Range("B57").Value = "spain"
$ws.Range("C57").Value = "laliga"
$ws.Range("D57").Value = "2023-2024"
$ws.Range("A57").Value = 56
$ws.Range("E57").Value = 45193.67708333334
$ws.Range("F57").Value = "Rayo Vallecano"
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = "Villarreal"
$ws.Range("I57").Value = 1
$ws.Range("J57").Value = 2.72
$ws.Range("K57").Value = "05/09/2023 12:02"
$ws.Range("L57").Value = 2.64
$ws.Range("M57").Value = "24/09/2023 16:11"
$ws.Range("N57").Value = 3.4
$ws.Range("O57").Value = "05/09/2023 12:02"
$ws.Range("P57").Value = 3.53
$ws.Range("Q57").Value = "24/09/2023 15:54"
$ws.Range("R57").Value = 2.69
$ws.Range("S57").Value = "05/09/2023 12:02"
$ws.Range("T57").Value = 2.78
$ws.Range("U57").Value = "24/09/2023 16:11"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/spain/laliga/rayo-vallecano-villarreal/6DkbGGJ4/"

# Row 58
$ws.Range("B58").Value = "spain"
$ws.Range("C58").Value = "laliga"
$ws.Range("D58").Value = "2023-2024"
$ws.Range("A58").Value = 57
$ws.Range("E58").Value = 45193.77083333334
$ws.Range("F58").Value = "Betis"
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = "Cadiz CF"
$ws.Range("I58").Value = 1
$ws.Range("J58").Value = 1.87
$ws.Range("K58").Value = "05/09/2023 12:02"
$ws.Range("L58").Value = 1.81
$ws.Range("M58").Value = "24/09/2023 18:26"
$ws.Range("N58").Value = 3.56
$ws.Range("O58").Value = "05/09/2023 12:02"
$ws.Range("P58").Value = 3.66
$ws.Range("Q58").Value = "24/09/2023 18:26"
$ws.Range("R58").Value = 4.54
$ws.Range("S58").Value = "05/09/2023 12:02"
$ws.Range("T58").Value = 5.07
$ws.Range("U58").Value = "24/09/2023 18:26"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/spain/laliga/betis-cadiz/IicoJIZo/"

# Row 59
$ws.Range("B59").Value = "spain"
$ws.Range("C59").Value = "laliga"
$ws.Range("D59").Value = "2023-2024"
$ws.Range("A59").Value = 58
$ws.Range("E59").Value = 45193.77083333334
$ws.Range("F59").Value = "Las Palmas"
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = "Granada CF"
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2.14
$ws.Range("K59").Value = "11/09/2023 13:19"
$ws.Range("L59").Value = 1.95
$ws.Range("M59").Value = "24/09/2023 18:19"
$ws.Range("N59").Value = 3.21
$ws.Range("O59").Value = "11/09/2023 13:19"
$ws.Range("P59").Value = 3.73
$ws.Range("Q59").Value = "24/09/2023 18:27"
$ws.Range("R59").Value = 3.61
$ws.Range("S59").Value = "11/09/2023 13:19"
$ws.Range("T59").Value = 4.1
$ws.Range("U59").Value = "24/09/2023 18:27"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/spain/laliga/las-palmas-granada-cf/tWsBDE3N/"

# Row 60
$ws.Range("B60").Value = "spain"
$ws.Range("C60").Value = "laliga"
$ws.Range("D60").Value = "2023-2024"
$ws.Range("A60").Value = 59
$ws.Range("E60").Value = 45193.875
$ws.Range("F60").Value = "Atl. Madrid"
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = "Real Madrid"
$ws.Range("I60").Value = 1
$ws.Range("J60").Value = 2.82
$ws.Range("K60").Value = "05/09/2023 12:02"
$ws.Range("L60").Value = 2.87
$ws.Range("M60").Value = "24/09/2023 20:55"
$ws.Range("N60").Value = 3.67
$ws.Range("O60").Value = "05/09/2023 12:02"
$ws.Range("P60").Value = 3.38
$ws.Range("Q60").Value = "24/09/2023 20:59"
$ws.Range("R60").Value = 2.33
$ws.Range("S60").Value = "05/09/2023 12:02"
$ws.Range("T60").Value = 2.65
$ws.Range("U60").Value = "24/09/2023 20:56"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/spain/laliga/atl-madrid-real-madrid/hCtzC1SI/"

